# Apply updated odds values to Sheet1 as per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 7).Value = 1.75  # G2: 1.73 -> 1.75
$ws.Cells.Item(2, 9).Value = 4.5  # I2: 4.75 -> 4.5
$ws.Cells.Item(2, 12).Value = 4.5  # L2: 4.75 -> 4.5
$ws.Cells.Item(2, 15).Value = 1.22  # O2: 1.2 -> 1.22
$ws.Cells.Item(2, 16).Value = 4.33  # P2: 4.5 -> 4.33
$ws.Cells.Item(2, 21).Value = 1.62  # U2: 1.67 -> 1.62
$ws.Cells.Item(2, 22).Value = 2.2  # V2: 2.1 -> 2.2
$ws.Cells.Item(2, 23).Value = 9  # W2: 8.5 -> 9
$ws.Cells.Item(2, 24).Value = 9.5  # X2: 9 -> 9.5
$ws.Cells.Item(2, 26).Value = 15  # Z2: 13 -> 15
$ws.Cells.Item(2, 35).Value = 23  # AI2: 26 -> 23
$ws.Cells.Item(2, 37).Value = 41  # AK2: 51 -> 41
$ws.Cells.Item(2, 40).Value = 4  # AN2: 3.75 -> 4
$ws.Cells.Item(2, 43).Value = 29  # AQ2: 26 -> 29
$ws.Cells.Item(2, 48).Value = 41  # AV2: 51 -> 41
$ws.Cells.Item(2, 49).Value = 451  # AW2: 501 -> 451
$ws.Cells.Item(2, 50).Value = 6  # AX2: 6.5 -> 6
$ws.Cells.Item(2, 51).Value = 21  # AY2: 23 -> 21
$ws.Cells.Item(2, 52).Value = 26  # AZ2: 29 -> 26
$ws.Cells.Item(2, 53).Value = 67  # BA2: 81 -> 67
# Row 3
$ws.Cells.Item(3, 10).Value = 2.5  # J3: 2.4 -> 2.5
$ws.Cells.Item(3, 12).Value = 5  # L3: 4.75 -> 5
$ws.Cells.Item(3, 15).Value = 1.33  # O3: 1.3 -> 1.33
$ws.Cells.Item(3, 16).Value = 3.4  # P3: 3.5 -> 3.4
$ws.Cells.Item(3, 17).Value = 2  # Q3: 2.03 -> 2
$ws.Cells.Item(3, 18).Value = 1.8  # R3: 1.87 -> 1.8
$ws.Cells.Item(3, 19).Value = 1.44  # S3: 1.4 -> 1.44
$ws.Cells.Item(3, 20).Value = 2.63  # T3: 2.75 -> 2.63
$ws.Cells.Item(3, 29).Value = 9  # AC3: 10 -> 9
$ws.Cells.Item(3, 30).Value = 6.5  # AD3: 7 -> 6.5
$ws.Cells.Item(3, 34).Value = 12  # AH3: 11 -> 12
$ws.Cells.Item(3, 38).Value = 41  # AL3: 34 -> 41
$ws.Cells.Item(3, 46).Value = 2.63  # AT3: 2.75 -> 2.63
# Row 4
$ws.Cells.Item(4, 7).Value = 2.35  # G4: 2.4 -> 2.35
$ws.Cells.Item(4, 9).Value = 3.5  # I4: 3.4 -> 3.5
$ws.Cells.Item(4, 10).Value = 3.2  # J4: 3.25 -> 3.2
$ws.Cells.Item(4, 26).Value = 21  # Z4: 23 -> 21
$ws.Cells.Item(4, 29).Value = 6.5  # AC4: 6 -> 6.5
$ws.Cells.Item(4, 30).Value = 6  # AD4: 5.5 -> 6
$ws.Cells.Item(4, 31).Value = 17  # AE4: 19 -> 17
$ws.Cells.Item(4, 34).Value = 8  # AH4: 7.5 -> 8
$ws.Cells.Item(4, 54).Value = 101  # BB4: 126 -> 101
$ws.Cells.Item(4, 56).Value = 126  # BD4: 151 -> 126
# Row 5
$ws.Cells.Item(5, 12).Value = 5  # L5: 4.75 -> 5
$ws.Cells.Item(5, 24).Value = 8  # X5: 8.5 -> 8
$ws.Cells.Item(5, 26).Value = 17  # Z5: 19 -> 17
$ws.Cells.Item(5, 34).Value = 9  # AH5: 8.5 -> 9
$ws.Cells.Item(5, 36).Value = 17  # AJ5: 15 -> 17
$ws.Cells.Item(5, 37).Value = 51  # AK5: 41 -> 51
# Row 6
$ws.Cells.Item(6, 13).Value = 1.1  # M6: 1.11 -> 1.1
$ws.Cells.Item(6, 14).Value = 7  # N6: 6.5 -> 7
$ws.Cells.Item(6, 15).Value = 1.44  # O6: 1.5 -> 1.44
$ws.Cells.Item(6, 16).Value = 2.63  # P6: 2.5 -> 2.63
$ws.Cells.Item(6, 56).Value = 126  # BD6: 151 -> 126
# Row 7
$ws.Cells.Item(7, 18).Value = 1.3  # R7: 1.33 -> 1.3
$ws.Cells.Item(7, 22).Value = 1.47  # V7: 1.5 -> 1.47
# Row 8
$ws.Cells.Item(8, 13).Value = 1.02  # M8: 1.03 -> 1.02
$ws.Cells.Item(8, 14).Value = 19  # N8: 15 -> 19
$ws.Cells.Item(8, 17).Value = 1.47  # Q8: 1.5 -> 1.47
$ws.Cells.Item(8, 21).Value = 1.87  # U8: 1.91 -> 1.87
$ws.Cells.Item(8, 22).Value = 1.77  # V8: 1.8 -> 1.77
# Row 11
$ws.Cells.Item(11, 13).Value = 1.05  # M11: 1.07 -> 1.05
$ws.Cells.Item(11, 14).Value = 11  # N11: 9 -> 11
# Row 12
$ws.Cells.Item(12, 17).Value = 2.15  # Q12: 2.1 -> 2.15
$ws.Cells.Item(12, 18).Value = 1.67  # R12: 1.7 -> 1.67
# Row 14
$ws.Cells.Item(14, 7).Value = 2.7  # G14: 2.6 -> 2.7
$ws.Cells.Item(14, 9).Value = 2.55  # I14: 2.63 -> 2.55
$ws.Cells.Item(14, 10).Value = 3.5  # J14: 3.4 -> 3.5
$ws.Cells.Item(14, 12).Value = 3.4  # L14: 3.5 -> 3.4
$ws.Cells.Item(14, 13).Value = 1.05  # M14: 1.07 -> 1.05
$ws.Cells.Item(14, 14).Value = 8.5  # N14: 9 -> 8.5
$ws.Cells.Item(14, 15).Value = 1.37  # O14: 1.4 -> 1.37
$ws.Cells.Item(14, 17).Value = 2.2  # Q14: 2.25 -> 2.2
$ws.Cells.Item(14, 18).Value = 1.65  # R14: 1.62 -> 1.65
$ws.Cells.Item(14, 24).Value = 13  # X14: 12 -> 13
$ws.Cells.Item(14, 25).Value = 11  # Y14: 10 -> 11
$ws.Cells.Item(14, 26).Value = 29  # Z14: 26 -> 29
$ws.Cells.Item(14, 29).Value = 8.5  # AC14: 8 -> 8.5
$ws.Cells.Item(14, 36).Value = 10  # AJ14: 11 -> 10
$ws.Cells.Item(14, 40).Value = 4.75  # AN14: 4.5 -> 4.75
$ws.Cells.Item(14, 41).Value = 17  # AO14: 15 -> 17
$ws.Cells.Item(14, 51).Value = 15  # AY14: 17 -> 15
$ws.Cells.Item(14, 52).Value = 26  # AZ14: 29 -> 26
$ws.Cells.Item(14, 55).Value = 201  # BC14: 251 -> 201
# Row 15
$ws.Cells.Item(15, 7).Value = 1.22  # G15: 1.2 -> 1.22
$ws.Cells.Item(15, 8).Value = 6  # H15: 6.25 -> 6
$ws.Cells.Item(15, 14).Value = 15  # N15: 19 -> 15
$ws.Cells.Item(15, 15).Value = 1.13  # O15: 1.17 -> 1.13
$ws.Cells.Item(15, 19).Value = 1.29  # S15: 1.25 -> 1.29
$ws.Cells.Item(15, 20).Value = 3.5  # T15: 3.75 -> 3.5
$ws.Cells.Item(15, 21).Value = 2.1  # U15: 2.2 -> 2.1
$ws.Cells.Item(15, 22).Value = 1.67  # V15: 1.62 -> 1.67
$ws.Cells.Item(15, 26).Value = 7.5  # Z15: 7 -> 7.5
$ws.Cells.Item(15, 32).Value = 67  # AF15: 81 -> 67
$ws.Cells.Item(15, 46).Value = 3.5  # AT15: 3.75 -> 3.5
$ws.Cells.Item(15, 47).Value = 10  # AU15: 11 -> 10
$ws.Cells.Item(15, 50).Value = 11  # AX15: 12 -> 11
$ws.Cells.Item(15, 52).Value = 41  # AZ15: 51 -> 41
$ws.Cells.Item(15, 53).Value = 251  # BA15: 301 -> 251
# Row 16
$ws.Cells.Item(16, 13).Value = 1.03  # M16: 1.05 -> 1.03
$ws.Cells.Item(16, 15).Value = 1.25  # O16: 1.29 -> 1.25
# Row 17
$ws.Cells.Item(17, 7).Value = 1.63  # G17: 1.65 -> 1.63
$ws.Cells.Item(17, 9).Value = 4.75  # I17: 5 -> 4.75
$ws.Cells.Item(17, 10).Value = 2.4  # J17: 2.38 -> 2.4
$ws.Cells.Item(17, 18).Value = 1.54  # R17: 1.57 -> 1.54
$ws.Cells.Item(17, 26).Value = 13  # Z17: 12 -> 13
$ws.Cells.Item(17, 34).Value = 10  # AH17: 11 -> 10
$ws.Cells.Item(17, 41).Value = 9.5  # AO17: 9 -> 9.5
$ws.Cells.Item(17, 42).Value = 26  # AP17: 23 -> 26
# Row 18
$ws.Cells.Item(18, 7).Value = 2.2  # G18: 2.05 -> 2.2
$ws.Cells.Item(18, 8).Value = 3  # H18: 3.1 -> 3
$ws.Cells.Item(18, 9).Value = 3.4  # I18: 4 -> 3.4
$ws.Cells.Item(18, 10).Value = 3.1  # J18: 2.88 -> 3.1
$ws.Cells.Item(18, 11).Value = 1.91  # K18: 1.95 -> 1.91
$ws.Cells.Item(18, 12).Value = 4.33  # L18: 4.75 -> 4.33
$ws.Cells.Item(18, 13).Value = 1.1  # M18: 1.11 -> 1.1
$ws.Cells.Item(18, 14).Value = 7  # N18: 6.5 -> 7
$ws.Cells.Item(18, 21).Value = 2.1  # U18: 2.2 -> 2.1
$ws.Cells.Item(18, 22).Value = 1.67  # V18: 1.62 -> 1.67
$ws.Cells.Item(18, 23).Value = 6  # W18: 5.5 -> 6
$ws.Cells.Item(18, 24).Value = 9.5  # X18: 8.5 -> 9.5
$ws.Cells.Item(18, 25).Value = 10  # Y18: 9.5 -> 10
$ws.Cells.Item(18, 26).Value = 21  # Z18: 17 -> 21
$ws.Cells.Item(18, 32).Value = 67  # AF18: 81 -> 67
$ws.Cells.Item(18, 34).Value = 8  # AH18: 8.5 -> 8
$ws.Cells.Item(18, 35).Value = 15  # AI18: 19 -> 15
$ws.Cells.Item(18, 36).Value = 13  # AJ18: 15 -> 13
$ws.Cells.Item(18, 38).Value = 34  # AL18: 41 -> 34
$ws.Cells.Item(18, 39).Value = 41  # AM18: 51 -> 41
$ws.Cells.Item(18, 40).Value = 4  # AN18: 3.75 -> 4
$ws.Cells.Item(18, 41).Value = 13  # AO18: 12 -> 13
$ws.Cells.Item(18, 43).Value = 51  # AQ18: 41 -> 51
$ws.Cells.Item(18, 47).Value = 9  # AU18: 9.5 -> 9
$ws.Cells.Item(18, 48).Value = 67  # AV18: 81 -> 67
$ws.Cells.Item(18, 50).Value = 5  # AX18: 5.5 -> 5
$ws.Cells.Item(18, 51).Value = 21  # AY18: 23 -> 21
$ws.Cells.Item(18, 52).Value = 34  # AZ18: 41 -> 34
$ws.Cells.Item(18, 53).Value = 67  # BA18: 81 -> 67
# Row 19
$ws.Cells.Item(19, 7).Value = 1.55  # G19: 1.57 -> 1.55
$ws.Cells.Item(19, 9).Value = 6  # I19: 5.75 -> 6
$ws.Cells.Item(19, 13).Value = 1.02  # M19: 1.03 -> 1.02
$ws.Cells.Item(19, 15).Value = 1.13  # O19: 1.18 -> 1.13
$ws.Cells.Item(19, 16).Value = 5  # P19: 4.5 -> 5
$ws.Cells.Item(19, 17).Value = 1.6  # Q19: 1.62 -> 1.6
$ws.Cells.Item(19, 18).Value = 2.3  # R19: 2.25 -> 2.3
$ws.Cells.Item(19, 50).Value = 7.5  # AX19: 7 -> 7.5
$ws.Cells.Item(19, 55).Value = 201  # BC19: 151 -> 201
